$d = $word.ActiveDocument

$replacements = @(
    @("374×9=", "896×9="),
    @("949×8=", "381×7="),
    @("533×5=", "563×4="),
    @("749×4=", "541×2="),
    @("740×7=", "867×2="),
    @("569×5=", "510×9="),
    @("940×6=", "754×3="),
    @("475×4=", "290×2="),
    @("108×2=", "286×4="),
    @("782×9=", "361×4="),
    @("488×8=", "228×5="),
    @("925×8=", "253×6="),
    @("550×7=", "211×5="),
    @("886×7=", "861×3="),
    @("829×3=", "589×5="),
    @("600×2=", "229×5="),
    @("450×2=", "436×4="),
    @("471×5=", "443×6="),
    @("202×3=", "829×4="),
    @("620×3=", "658×2="),
    @("980×7=", "883×5="),
    @("266×8=", "464×5="),
    @("934×9=", "110×6="),
    @("958×4=", "596×3="),
    @("667×5=", "897×9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
